# ---------------------------------------------------------------------------
# Update racial-misclassification language to neutral "50M voters" wording
# and relocate the "Field Director - The Feldman Group" role to sit right
# after the "Research Director - PCCC" role (immediately before
# "Software Engineer - Salsa Labs").
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the 1-based Paragraphs index whose Range contains a given
# document character position.
# ---------------------------------------------------------------------------
function Get-ParagraphIndexAt($pos) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($pos -ge $p.Range.Start -and $pos -lt $p.Range.End) {
            return $i
        }
    }
    return $count
}

# ---------------------------------------------------------------------------
# Change 1: Professional summary paragraph — plain text swap.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "affecting all Black and Asian-American voters, developed geospatial ML",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "affecting 50M voters, developed geospatial ML", 2
) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: "Impact: Corrected demographic data ..." project line — plain
# text swap.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Impact: Corrected demographic data affecting all Black and Asian-American voters, improved",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Impact: Corrected demographic data affecting 50M voters nationwide, improved", 2
) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: Siege Analytics bullet point — split the run so that "50M"
# becomes its own bold, colored run (matching the styling of the other
# statistic call-outs in this bullet).
# ---------------------------------------------------------------------------
$bulletTarget = $d.Content.Duplicate
$bulletTarget.Find.Execute(
    "all Black and Asian-American",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
) | Out-Null

$bulletRange = $d.Range($bulletTarget.Start, $bulletTarget.End)
$bulletRange.Text = "50M"

$bulletBold = $d.Range($bulletTarget.Start, $bulletTarget.Start + 3)
$bulletBold.Bold = 1
$bulletBold.Font.Color = 5258796   # RGB(0x2C, 0x3E, 0x50) == w:color val="2C3E50"

# ---------------------------------------------------------------------------
# Change 4: Move the "Field Director - The Feldman Group" entry (heading +
# its 4 paragraphs) from its old spot (right before "KEY PROJECTS") to a new
# spot right before "Software Engineer - Salsa Labs".
# ---------------------------------------------------------------------------

# --- Locate the section to move (source) ---
$srcStartFind = $d.Content.Duplicate
$srcStartFind.Find.Execute(
    "Field Director - The Feldman Group (Austin, TX) | 2011 - 2012",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
) | Out-Null

$srcEndFind = $d.Content.Duplicate
$srcEndFind.Find.Execute(
    "Created custom reports and data visualizations based on specific client requirements",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
) | Out-Null
$srcEndFind.MoveEnd(1, 1)   # include the trailing paragraph mark

$srcRange = $d.Range($srcStartFind.Start, $srcEndFind.End)
$srcLength = $srcRange.End - $srcRange.Start

# --- Locate the destination (right before the Salsa Labs heading) ---
$destFind = $d.Content.Duplicate
$destFind.Find.Execute(
    "Software Engineer - Salsa Labs (Washington, DC) | January 2011 - August 2011",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
) | Out-Null
$destPos = $destFind.Start

# --- Copy the formatted section to the new location ---
$sectionFormattedText = $srcRange.FormattedText
$insertionPoint = $d.Range($destPos, $destPos)
$insertionPoint.FormattedText = $sectionFormattedText

# The copy above drops the paragraph style of the heading line (FormattedText
# only preserves character-level formatting), so re-apply it explicitly.
$newHeadingIndex = Get-ParagraphIndexAt($destPos)
$newHeadingParagraph = $d.Paragraphs.Item($newHeadingIndex)
$newHeadingParagraph.Style = "Heading 3"

# --- Remove the original section, which has shifted forward by the length
#     of the copy we just inserted ahead of it. Re-find it defensively
#     rather than trusting raw arithmetic. ---
$searchAfter = $destPos + $srcLength

$origStartFind = $d.Content.Duplicate
$origStartFind.Start = $searchAfter
$origStartFind.Find.Execute(
    "Field Director - The Feldman Group (Austin, TX) | 2011 - 2012",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
) | Out-Null

$origEndFind = $d.Content.Duplicate
$origEndFind.Start = $searchAfter
$origEndFind.Find.Execute(
    "Created custom reports and data visualizations based on specific client requirements",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
) | Out-Null
$origEndFind.MoveEnd(1, 1)

$origRange = $d.Range($origStartFind.Start, $origEndFind.End)
$origRange.Delete()

Write-Output "Edit complete. Paragraph count: $($d.Paragraphs.Count)"
